$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, $row, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $i + 1).Value2 = $values[$i]
    }
}

# ---------------------------------------------------------------
# Sheet "Home win": 2 -> 4 data rows
# ---------------------------------------------------------------
$wsHome = $wb.Worksheets.Item("Home win")
Set-Row $wsHome 2 @("04-02-2025 00:30","COLOMBIA","PRIMERA A","Deportivo Pereira - Alianza Petrolera",70,1.86)
Set-Row $wsHome 3 @("04-02-2025 20:45","ENGLAND","LEAGUE TWO","Salford City - Bromley",73.3,2)
Set-Row $wsHome 4 @("04-02-2025 20:45","SCOTLAND","LEAGUE TWO","Elgin City - Bonnyrigg Rose Athletic",73.3,2)
Set-Row $wsHome 5 @("04-02-2025 20:00","ENGLAND","PREMIER LEAGUE CUP","Ipswich Town U21 - Watford U21",73.3,1.73)

# ---------------------------------------------------------------
# Sheet "Draw": 2 data rows, values refreshed
# ---------------------------------------------------------------
$wsDraw = $wb.Worksheets.Item("Draw")
Set-Row $wsDraw 2 @("04-02-2025 20:45","ENGLAND","CHAMPIONSHIP","Burnley - Oxford United",70,3.9)
Set-Row $wsDraw 3 @("04-02-2025 20:45","ENGLAND","NON LEAGUE PREMIER - NORTHERN","Prescot Cables - Rylands",60,3.3)

# ---------------------------------------------------------------
# Sheet "Btts": 5 -> 4 data rows
# ---------------------------------------------------------------
$wsBtts = $wb.Worksheets.Item("Btts")
Set-Row $wsBtts 2 @("04-02-2025 23:00","BRAZIL","CEARENSE - 1","Horizonte - Pague Menos",80,1.85)
Set-Row $wsBtts 3 @("05-02-2025 00:00","CHILE","COPA CHILE","Deportes Limache - Union San Felipe",83.3,1.73)
Set-Row $wsBtts 4 @("04-02-2025 21:10","FRANCE","COUPE DE FRANCE","Le Mans - Paris Saint Germain",90,2.1)
Set-Row $wsBtts 5 @("04-02-2025 19:00","FRANCE","COUPE DE FRANCE","Lille - Dunkerque",78.3,1.85)
$wsBtts.Range("A6:F6").ClearContents()

# ---------------------------------------------------------------
# Sheet "Over_Under": 6 -> 12 data rows
# ---------------------------------------------------------------
$wsOU = $wb.Worksheets.Item("Over_Under")
Set-Row $wsOU 2  @("04-02-2025 20:45","ENGLAND","EFL TROPHY","Stevenage - Birmingham",80,2,13.3,3.4)
Set-Row $wsOU 3  @("04-02-2025 20:45","ENGLAND","NATIONAL LEAGUE - NORTH","Chorley - Buxton",65,1.95,50,3.3)
Set-Row $wsOU 4  @("04-02-2025 20:45","ENGLAND","NATIONAL LEAGUE - NORTH","Scarborough Athletic - Radcliffe",70,1.65,60,2.6)
Set-Row $wsOU 5  @("04-02-2025 20:45","ENGLAND","NON LEAGUE PREMIER - ISTHMIAN","Cray Valley PM - Lewes",73.3,1.57,60,2.4)
Set-Row $wsOU 6  @("04-02-2025 20:45","ENGLAND","NON LEAGUE PREMIER - SOUTHERN SOUTH","Merthyr Town - Basingstoke Town",55,1.45,55,2.25)
Set-Row $wsOU 7  @("04-02-2025 19:00","FRANCE","COUPE DE FRANCE","Lille - Dunkerque",70,1.8,60,3)
Set-Row $wsOU 8  @("04-02-2025 18:00","ISRAEL","LIGA LEUMIT","Hapoel Ramat Gan - Bnei Yehuda",80,1.75,55,2.88)
Set-Row $wsOU 9  @("04-02-2025 18:00","ISRAEL","LIGA LEUMIT","Hapoel Afula - Hapoel Rishon LeZion",70,2,60,3.5)
Set-Row $wsOU 10 @("04-02-2025 18:00","ISRAEL","LIGA LEUMIT","Hapoel Kfar Saba - Hapoel Ramat HaSharon",65,1.73,60,2.8)
Set-Row $wsOU 11 @("04-02-2025 15:00","PORTUGAL","LIGA REVELAÇÃO U23","Benfica U23 - Torreense U23",50,1.65,50,2.6)
Set-Row $wsOU 12 @("04-02-2025 20:30","SWITZERLAND","SUPER LEAGUE","FC Sion - Servette FC",70,1.83,55,3.1)
Set-Row $wsOU 13 @("04-02-2025 10:30","WORLD","FRIENDLIES CLUBS","Slavia Praha II - Příbram",80,1.57,53.3,2.4)

Write-Output "edits applied"
